$d = $word.ActiveDocument

# --- Table resize: tblInd 0 -> 113 (dxa) and tblGrid column widths ---
$t = $d.Tables.Item(1)

# w:tblInd w:w="113" w:type="dxa"  (COM LeftIndent is in points; 20 twips = 1 pt)
$t.Rows.LeftIndent = 113 / 20.0

# w:gridCol widths (dxa -> points)
$t.Columns.Item(1).Width = 4232 / 20.0
$t.Columns.Item(2).Width = 3391 / 20.0
$t.Columns.Item(3).Width = 3391 / 20.0

# --- Section / page setup: nudge <w:cols> towards the canonical single-column
# self-closed form (adds w:space="0", matching the target) ---
$sec = $d.Sections.Item(1)
$ps = $sec.PageSetup
$ps.TextColumns.Spacing = 0

Write-Output "table indent/grid + section columns updated"
